$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from the existing column G cell (G4) to the new cells
# so the new date cells use the same style index instead of creating a
# duplicate numFmt entry.
$ws.Range("G4").Copy()
$ws.Range("G5:G6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 5
$ws.Range("A5").Value = 9846.56
$ws.Range("B5").Value = 10136.459999999999
$ws.Range("C5").Value = 113.86
$ws.Range("D5").Value = 110.6
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -2.86
$ws.Range("G5").Value = 42607.884247685186
$ws.Range("H5").Value = $false

# Row 6
$ws.Range("A6").Value = 9592.52
$ws.Range("B6").Value = 9846.56
$ws.Range("C6").Value = 110.77
$ws.Range("D6").Value = 107.91
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -2.58
$ws.Range("G6").Value = 42608.616331018522
$ws.Range("H6").Value = $false
